$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = 257
$ws.Range("C8").Value = 338
$ws.Range("C9").Value = 438
$ws.Range("C10").Value = 521
$ws.Range("C12").Value = 696
$ws.Range("C13").Value = 781
$ws.Range("C14").Value = 893
$ws.Range("I14").Value = 3109
$ws.Range("C15").Value = 886
$ws.Range("C17").Value = 1008
$ws.Range("I17").Value = 3907
$ws.Range("C18").Value = 1115
$ws.Range("C19").Value = 1222
$ws.Range("C20").Value = 1303
$ws.Range("I20").Value = 5629
$ws.Range("C21").Value = 1490
$ws.Range("I21").Value = 6910
$ws.Range("I22").Value = 7648
$ws.Range("C23").Value = 1959
$ws.Range("I23").Value = 8912
$ws.Range("C24").Value = 2234
$ws.Range("I24").Value = 10383
$ws.Range("C25").Value = 2501
$ws.Range("I25").Value = 10866
$ws.Range("C26").Value = 2760
$ws.Range("I26").Value = 13875
$ws.Range("C27").Value = 3180
$ws.Range("I27").Value = 16381
$ws.Range("C28").Value = 3624
$ws.Range("I28").Value = 18175
$ws.Range("C29").Value = 4186
$ws.Range("I29").Value = 19563
$ws.Range("C30").Value = 4751
$ws.Range("I30").Value = 21816
$ws.Range("C31").Value = 5118
$ws.Range("I31").Value = 24856
$ws.Range("C32").Value = 5618
$ws.Range("I32").Value = 28297
$ws.Range("C33").Value = 5953
$ws.Range("C34").Value = 6186
$ws.Range("I34").Value = 32389
$ws.Range("C35").Value = 6451
$ws.Range("I35").Value = 35499
$ws.Range("C36").Value = 7392
$ws.Range("I36").Value = 42616
$ws.Range("C37").Value = 11880
$ws.Range("I37").Value = 67340
